$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "mfd_hab2" values ("Fjords") from rows 2-11 in column O,
# leaving the header (O1) and the rest of the data untouched.
$ws.Range("O2:O11").ClearContents()
